# "Add files via upload" — append the two new Facebook-post log rows that
# were added at the bottom of the tracking sheet (rows 99 and 100), and
# leave the cursor/selection where the author left it (cell E102) after
# scrolling the frozen-header view down to row 77.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 99 : 2021-02-15 09:13 AM -------------------------------------
$ws.Cells.Item(99, 1).Value  = 44242                        # Date
$ws.Cells.Item(99, 2).Value  = 0.3840277777777778           # Time
$ws.Cells.Item(99, 3).Value  = "Friends"                    # Privacy
$ws.Cells.Item(99, 4).Value  = "DENSE FOG ☁️(Also at Dentist)"           # Post
$ws.Cells.Item(99, 5).Value  = "10108072387916819"          # Reference
$ws.Cells.Item(99, 6).Value  = 6                             # Like
$ws.Cells.Item(99, 7).Value  = 0                             # Love
$ws.Cells.Item(99, 8).Value  = 0                             # Wow
$ws.Cells.Item(99, 9).Value  = 0                             # Haha
$ws.Cells.Item(99, 10).Value = 0                             # Sad
$ws.Cells.Item(99, 11).Value = 0                             # Care
$ws.Cells.Item(99, 12).Value = 1                             # Angry
$ws.Cells.Item(99, 13).Value = 6                             # Comments

# --- Row 100 : 2021-02-15 08:12 PM ------------------------------------
$ws.Cells.Item(100, 1).Value  = 44242
$ws.Cells.Item(100, 2).Value  = 0.84166666666666667
$ws.Cells.Item(100, 3).Value  = "Friends"
$ws.Cells.Item(100, 4).Value  = "#email 📬⚡️📭 :: #doubleexposure #densefog #ncwinter #drizzle #night #lights"
$ws.Cells.Item(100, 5).Value  = "10108073812921099"
$ws.Cells.Item(100, 6).Value  = 7
$ws.Cells.Item(100, 7).Value  = 0
$ws.Cells.Item(100, 8).Value  = 1
$ws.Cells.Item(100, 9).Value  = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = 1

# Leave the view/selection the way the workbook was saved: scrolled down
# with E102 (the row right after the new data) selected.
$ws.Range("E102").Select()
